# Report alles auf Englisch - add new translation rows for
# PT/PD, Summe über alle Projekte/Sum over all projects,
# Platzhalter/placeholder, zugeordnet/assigned
# Column A = German text, Column B = English text (matches existing layout)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A274").Value = "PT"
$ws.Range("B274").Value = "PD"

$ws.Range("A275").Value = "Summe über alle Projekte"
$ws.Range("B275").Value = "Sum over all projects"

$ws.Range("A276").Value = "Platzhalter"
$ws.Range("B276").Value = "placeholder"

$ws.Range("A277").Value = "zugeordnet"
$ws.Range("B277").Value = "assigned"

# Update the view state to match the new selection / scroll position
$ws.Range("B274").Select()
